$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 70, shifting existing rows 70-73 down to 71-74.
$ws.Rows.Item(70).EntireRow.Insert()

# Populate the newly inserted row 70 with the new weekly record.
$ws.Range("A70").Value = 10
$ws.Range("B70").Value = "Vega Modelo de Temuco"
$ws.Range("C70").Value = "La Araucanía"
$ws.Range("D70").Value = 45041
$ws.Range("E70").Value = 9
$ws.Range("F70").Value = 100112042
$ws.Range("G70").Value = "Locoto"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 150
$ws.Range("K70").Value = 4400
$ws.Range("L70").Value = 4400
$ws.Range("M70").Value = 4400
$ws.Range("N70").Value = "$/kilo"
$ws.Range("O70").Value = "Región de Arica y Parinacota"
$ws.Range("P70").Value = 4400
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Range("D70").NumberFormat = $ws.Range("D69").NumberFormat
